$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2 (description): add two new helper columns F/G (blank centered
# merged header + two note cells reusing the existing "Neutral"/"Bad" look) ---
$ws2.Range("F1:G1").ColumnWidth = 22.6
$ws2.Range("F1:G1").HorizontalAlignment = -4108
$ws2.Range("F1:G1").Merge()

$ws2.Range("A2").Copy()
$ws2.Range("F2").PasteSpecial(-4122)
$ws2.Range("F2").Value = "เป็นค่าว่างได้"

$ws2.Range("B2").Copy()
$ws2.Range("G2").PasteSpecial(-4122)
$ws2.Range("G2").Value = "เพิ่มข้อมูล ต้องไม่เป็นค่าว่าง"

# --- Sheet2: update header text for the "เลขที่ (null)" column (done last so
# the new shared string lands after the other newly-added ones) ---
$ws2.Range("A1").Value = "เลขที่  (null)  ตัวเลขเท่านั้น"

# --- Sheet2: move the selection cursor ---
$ws2.Range("B6").Select()

# --- Switch the active sheet: Sheet1 becomes the active tab/selection ---
$ws1.Activate()
$ws1.Range("A7").Select()
